$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '63.155.31'
$ws.Cells.Item(2, 5).Value = '  -7.87%  '
Set-TextValue 3 4 '3.518.48'
$ws.Cells.Item(3, 5).Value = '  -3.45%  '
$ws.Cells.Item(4, 5).Value = '  +0.33%  '
Set-TextValue 5 4 '390.28'
$ws.Cells.Item(5, 5).Value = '  -7.07%  '
Set-TextValue 6 4 '122.02'
$ws.Cells.Item(6, 5).Value = '  -7.64%  '
Set-TextValue 7 4 '3.509.76'
$ws.Cells.Item(7, 5).Value = '  -3.22%  '
Set-TextValue 8 4 '0.585'
$ws.Cells.Item(8, 5).Value = '  -11.86%  '
$ws.Cells.Item(9, 5).Value = '  +0.08%  '
Set-TextValue 10 4 '0.677'
$ws.Cells.Item(10, 5).Value = '  -12.78%  '
Set-TextValue 11 4 '0.149'
$ws.Cells.Item(11, 5).Value = '  -26.40%  '
Set-TextValue 12 4 '0.0000317'
$ws.Cells.Item(12, 5).Value = '  -28.49%  '
Set-TextValue 13 4 '38.59'
$ws.Cells.Item(13, 5).Value = '  -9.08%  '
Set-TextValue 14 4 '4.094.54'
$ws.Cells.Item(14, 5).Value = '  -3.04%  '
Set-TextValue 15 4 '9.12'
$ws.Cells.Item(15, 5).Value = '  -8.28%  '
$ws.Cells.Item(16, 5).Value = '  -2.95%  '
Set-TextValue 17 4 '3.504.23'
$ws.Cells.Item(17, 5).Value = '  -4.07%  '
Set-TextValue 18 4 '12.79'
$ws.Cells.Item(18, 5).Value = '  +2.14%  '
Set-TextValue 19 4 '18.66'
$ws.Cells.Item(19, 5).Value = '  -7.61%  '
Set-TextValue 20 4 '63.315.24'
$ws.Cells.Item(20, 5).Value = '  -7.63%  '
Set-TextValue 21 4 '1.01'
$ws.Cells.Item(21, 5).Value = '  -11.00%  '
Set-TextValue 22 4 '389.31'
$ws.Cells.Item(22, 5).Value = '  -15.88%  '
Set-TextValue 23 4 '13.81'
$ws.Cells.Item(23, 5).Value = '  +2.93%  '
Set-TextValue 24 4 '80.31'
$ws.Cells.Item(24, 5).Value = '  -11.66%  '
Set-TextValue 25 4 '2.87'
$ws.Cells.Item(25, 5).Value = '  -7.11%  '
Set-TextValue 26 4 '5.46'
$ws.Cells.Item(26, 5).Value = '  +10.99%  '
Set-TextValue 27 4 '33.52'
$ws.Cells.Item(27, 5).Value = '  -6.70%  '
Set-TextValue 28 4 '2.97'
$ws.Cells.Item(28, 5).Value = '  -10.74%  '
Set-TextValue 29 4 '8.67'
$ws.Cells.Item(29, 5).Value = '  -15.74%  '
Set-TextValue 30 4 '2.64'
$ws.Cells.Item(30, 5).Value = '  -5.39%  '
Set-TextValue 31 4 '11.75'
$ws.Cells.Item(31, 5).Value = '  -5.19%  '
$ws.Cells.Item(32, 5).Value = '  -7.13%  '
Set-TextValue 33 4 '6.78'
$ws.Cells.Item(33, 5).Value = '  -6.44%  '
Set-TextValue 34 4 '0.149'
$ws.Cells.Item(34, 5).Value = '  -5.87%  '
Set-TextValue 35 4 '0.999'
$ws.Cells.Item(35, 5).Value = '  +0.06%  '
Set-TextValue 36 4 '36.73'
$ws.Cells.Item(36, 5).Value = '  -9.34%  '
Set-TextValue 37 4 '53.60'
$ws.Cells.Item(37, 5).Value = '  -4.46%  '
Set-TextValue 38 4 '0.0436'
$ws.Cells.Item(38, 5).Value = '  -10.55%  '
Set-TextValue 39 4 '1.00'
$ws.Cells.Item(39, 5).Value = '  +0.08%  '
Set-TextValue 40 4 '2.67'
$ws.Cells.Item(40, 5).Value = '  +5.19%  '
Set-TextValue 41 4 '0.130'
$ws.Cells.Item(41, 5).Value = '  -13.04%  '
Set-TextValue 42 4 '141.74'
$ws.Cells.Item(42, 5).Value = '  -4.87%  '
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 43 4 '26.01'
$ws.Cells.Item(43, 5).Value = '  +19.97%  '
$ws.Cells.Item(44, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 44 4 '3.04'
$ws.Cells.Item(44, 5).Value = '  +14.35%  '
$ws.Cells.Item(45, 5).Value = '  -28.53%  '
Set-TextValue 46 4 '2.49'
$ws.Cells.Item(46, 5).Value = '  -10.13%  '
Set-TextValue 47 4 '1.95'
$ws.Cells.Item(47, 5).Value = '  -0.82%  '
Set-TextValue 48 4 '3.07'
$ws.Cells.Item(48, 5).Value = '  -6.43%  '
$ws.Cells.Item(49, 5).Value = '  -5.23%  '
Set-TextValue 50 4 '2.65'
$ws.Cells.Item(50, 5).Value = '  -10.58%  '
$ws.Cells.Item(51, 5).Value = '  -10.29%  '
